$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) cells to retain text type (not auto-convert to numbers),
# matching the source workbook where prices are stored as inline strings.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "48.210.58"
$ws.Range("E2").Value = "  +2.26%  "
$ws.Range("D3").Value = "2.522.43"
$ws.Range("E3").Value = "  +1.37%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("B5").Value = "Solana"
$ws.Range("C5").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D5").Value = "109.96"
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "323.27"
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("D7").Value = "0.534"
$ws.Range("E7").Value = "  +2.07%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "0.556"
$ws.Range("E9").Value = "  +4.05%  "
$ws.Range("D10").Value = "40.72"
$ws.Range("E10").Value = "  +5.12%  "
$ws.Range("D11").Value = "20.42"
$ws.Range("E11").Value = "  +12.06%  "
$ws.Range("D12").Value = "0.0826"
$ws.Range("E12").Value = "  +1.91%  "
$ws.Range("E13").Value = "  +1.23%  "
$ws.Range("E14").Value = "  +1.81%  "
$ws.Range("D15").Value = "2.918.63"
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("D16").Value = "2.521.61"
$ws.Range("E16").Value = "  +1.32%  "
$ws.Range("D17").Value = "0.855"
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("D18").Value = "48.042.66"
$ws.Range("E18").Value = "  +2.06%  "
$ws.Range("D19").Value = "13.23"
$ws.Range("E19").Value = "  +4.19%  "
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("E21").Value = "  +1.32%  "
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("D23").Value = "72.13"
$ws.Range("E23").Value = "  +2.29%  "
$ws.Range("D24").Value = "263.87"
$ws.Range("E24").Value = "  +7.29%  "
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("D26").Value = "26.18"
$ws.Range("E26").Value = "  +1.83%  "
$ws.Range("E27").Value = "  -0.28%  "
$ws.Range("E28").Value = "  +1.02%  "
$ws.Range("E29").Value = "  +2.89%  "
$ws.Range("E30").Value = "  -2.79%  "
$ws.Range("D31").Value = "36.44"
$ws.Range("E31").Value = "  +3.62%  "
$ws.Range("D32").Value = "49.67"
$ws.Range("E32").Value = "  -0.45%  "
$ws.Range("D33").Value = "19.94"
$ws.Range("E33").Value = "  -0.58%  "
$ws.Range("D34").Value = "5.39"
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D36").Value = "0.0793"
$ws.Range("E36").Value = "  +1.08%  "
$ws.Range("D37").Value = "1.98"
$ws.Range("E37").Value = "  +1.43%  "
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("E39").Value = "  +1.54%  "
$ws.Range("E40").Value = "  +0.87%  "
$ws.Range("D41").Value = "120.88"
$ws.Range("E41").Value = "  +1.39%  "
$ws.Range("D42").Value = "22.07"
$ws.Range("E42").Value = "  +3.90%  "
$ws.Range("E43").Value = "  -1.29%  "
$ws.Range("E44").Value = "  +2.05%  "
$ws.Range("D45").Value = "2.018.79"
$ws.Range("E45").Value = "  +1.78%  "
$ws.Range("E46").Value = "  +4.85%  "
$ws.Range("D47").Value = "1.91"
$ws.Range("E47").Value = "  +6.80%  "
$ws.Range("E48").Value = "  +0.58%  "
$ws.Range("D49").Value = "9.10"
$ws.Range("D50").Value = "5.26"
$ws.Range("E50").Value = "  +2.78%  "
$ws.Range("D51").Value = "79.24"
$ws.Range("E51").Value = "  +2.81%  "

# Restore default (unstyled) appearance now that text values are locked in.
$ws.Range("D2:D51").Style = "Normal"

